# Commit: "Added errors and warning messages to test sheets"
#
# This adds two new worksheets ("Errors" and "Warnings") after the existing
# "Classes" sheet. "Errors" gets four header-validation messages (stored as
# shared strings, appended after the existing 7), and "Warnings" is left
# empty but becomes the active/selected sheet, matching the workbook's
# previous "Classes" tab-selection having moved off of "Classes".

$wb = $excel.ActiveWorkbook

# --- Try to restore the author's last-saved window geometry (best effort;
#     some hosts don't persist ActiveWindow geometry back into bookViews,
#     but it is harmless to set).
$win = $excel.ActiveWindow
$win.Left = 28740
$win.Top = -20060
$win.Width = 30500
$win.Height = 18760

# --- Add "Errors" sheet right after "Classes" ---------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$errorsSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$errorsSheet.Name = "Errors"

# Four header-validation error messages. Each literally begins with a
# single quote character (not an Excel "text prefix" marker) followed by
# the message text and a trailing comma, e.g. PHP var_export-style output.
# Doubling the leading quote on input makes Excel store one literal quote
# character, then resetting the cell Style strips the auto-applied
# "quote prefix" number format so the cell keeps the default style.
$errorsSheet.Range("A1").Value = "''Sheet ""Classes"" Row: 1 Column ""A"" in the header is not labeled as ""DDBNNN""',"
$errorsSheet.Range("A1").Style = "Normal"
$errorsSheet.Range("A2").Value = "''Sheet ""Classes"" Row: 1 Column ""B"" in the header is not labeled as ""TITLE""',"
$errorsSheet.Range("A2").Style = "Normal"
$errorsSheet.Range("A3").Value = "''Sheet ""Classes"" Row: 1 Column ""C"" in the header is not labeled as ""OFF CLS""',"
$errorsSheet.Range("A3").Style = "Normal"
$errorsSheet.Range("A4").Value = "''Sheet ""Classes"" Row: 1 Column ""D"" in the header is not labeled as ""SUB CLASSES""',"
$errorsSheet.Range("A4").Style = "Normal"

$errorsSheet.Range("A11").Select() | Out-Null

# --- Add "Warnings" sheet right after "Errors" ---------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$warningsSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$warningsSheet.Name = "Warnings"

# Left empty (no warnings), but it is the sheet that was active/selected
# when the workbook was last saved.
$warningsSheet.Range("D43").Select() | Out-Null

Write-Host "Sheets now:" ($wb.Worksheets | ForEach-Object { $_.Name }) -join ", "
